$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 held only the "5840671 - Francisco José Moreira Chaves" value
# (columns B/C, no label in A). That row is removed entirely; everything below
# shifts up by one row.
$ws.Rows(13).Delete() | Out-Null

# After the shift, patch the individual label values that changed content.

# "Objetivos:" (row 10) now shows the docent's name instead of the long
# objectives paragraph.
$ws.Range("B10").Value2 = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C10").Value2 = "5840671 - Francisco José Moreira Chaves"

# "Programa resumido:" (row 13) now just says "Semestral".
$ws.Range("B13").Value2 = "Semestral"
$ws.Range("C13").Value2 = "Semestral"

# "Programa:" (row 15) now shows a date instead of the long syllabus text.
$ws.Range("B15").Value2 = "01/01/2018"
$ws.Range("C15").Value2 = "01/01/2018"

# "Método:" (row 18) now shows the docent's name.
$ws.Range("B18").Value2 = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C18").Value2 = "5840671 - Francisco José Moreira Chaves"

# "Critério:" (row 19) now shows the old "Método:" text.
$ws.Range("B19").Value2 = "Por meio de aulas presenciais, com apresentação dos fundamentos e exemplos ou casos, e também pela apresentação de trabalhos em equipes.Justificativa: adequação do método de avaliação."
$ws.Range("C19").Value2 = "Por meio de aulas presenciais, com apresentação dos fundamentos e exemplos ou casos, e também pela apresentação de trabalhos em equipes.Justificativa: adequação do método de avaliação."

# "Norma de recuperação:" (row 20) now shows the old "Critério:" text.
$ws.Range("B20").Value2 = "A Avaliação será: MF = (P1 + P2)/2; Onde: P1: Trabalho; P2: Trabalho. Poderá haver também prova individual sobre os fundamentos."
$ws.Range("C20").Value2 = "A Avaliação será: MF = (P1 + P2)/2; Onde: P1: Trabalho; P2: Trabalho. Poderá haver também prova individual sobre os fundamentos."

# "Bibliografia:" (row 21) now shows the old "Norma de recuperação:" text
# ("Prova de exame."); the long bibliography list is dropped entirely.
$ws.Range("B21").Value2 = "Prova de exame."
$ws.Range("C21").Value2 = "Prova de exame."
